$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Home_Page): clear Runmode-status column B
$ws.Range("B4").Value = ""

# Row 5 (MyProfile): clear B
$ws.Range("B5").Value = ""

# Row 6 (Dashboard_Page): C6 Y -> N
$ws.Range("C6").Value = "N"

# Row 7 (RightHandPanel): clear B
$ws.Range("B7").Value = ""

# Row 9 (MyVehicles): clear B, C9 N -> Y
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "Y"

# Row 10 (MyMojios): clear B
$ws.Range("B10").Value = ""

# Row 11 (ClaimMojio): clear B
$ws.Range("B11").Value = ""

# Row 14 (MojioShop): change description text
$ws.Range("B14").Value = "Mojio Store-- not a part of My Mojio"

# Row 15 (MyOrders): change description text
$ws.Range("B15").Value = "My orders-- not a part of My Mojio"

# Row 16 (LocateMojio -> LocateVehicle): rename TSID, change description
$ws.Range("A16").Value = "LocateVehicle"
$ws.Range("B16").Value = "In Progress"

# Row 18 (SendFeedback): clear B
$ws.Range("B18").Value = ""

# Row 19 (Settings): clear B
$ws.Range("B19").Value = ""

# Row 20 (AdminSearch): clear B
$ws.Range("B20").Value = ""

# Row 21 (AdminDashboard): clear B
$ws.Range("B21").Value = ""

# Row 24 (ImportSIMs): clear B
$ws.Range("B24").Value = ""

# Row 25 (ImportMojio): clear B
$ws.Range("B25").Value = ""

# Row 26 (ExportEvents): clear B
$ws.Range("B26").Value = ""

# Update active cell selection to C9
$ws.Range("C9").Select()
